$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.207.15"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "2.425.37"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'554.31"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'137.30"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("E9").Value = "  -1.20%  "

$ws.Range("D10").Value = "'5.73"

$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("E12").Value = "  -1.27%  "

$ws.Range("D13").Value = "'24.93"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "2.857.76"
$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "60.122.09"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "'0.0000139"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "2.408.44"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").Value = "'11.28"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "'4.51"
$ws.Range("E19").Value = "  +2.91%  "

$ws.Range("D20").Value = "'327.51"
$ws.Range("E20").Value = "  -1.68%  "

$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "'65.32"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'0.178"
$ws.Range("E24").Value = "  +4.33%  "

$ws.Range("D25").Value = "'8.67"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  +5.52%  "

$ws.Range("D28").Value = "0.0₃0776"
$ws.Range("E28").Value = "  -1.00%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "'170.49"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("D31").Value = "'6.12"
$ws.Range("E31").Value = "  -2.31%  "

$ws.Range("D32").Value = "'0.404"
$ws.Range("E32").Value = "  -3.17%  "

$ws.Range("D33").Value = "'1.07"
$ws.Range("E33").Value = "  +2.23%  "

$ws.Range("D34").Value = "'18.56"
$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("E35").Value = "  +3.08%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "'326.76"
$ws.Range("E39").Value = "  +4.29%  "

$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").Value = "'145.48"
$ws.Range("E41").Value = "  +4.53%  "

$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.0964"
$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'19.89"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "'0.577"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").Value = "'0.0224"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "'11.05"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("E51").Value = "  -0.66%  "
